$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "06/11/2025"
$ws.Cells.Item(22, 1).Style = "Normal"
$ws.Cells.Item(22, 2).Value = 0.0004524999999999998
$ws.Cells.Item(22, 3).Value = 109392.2651933702
$ws.Cells.Item(22, 4).Value = 49.5
